# Generate Report for Handback
#
# A handback transform failed for the 5e890ed8 file: the zip the
# localization vendor returned contained a member (yyhf4t35.vza) whose
# name didn't match the handoff file it was supposed to correspond to.
# Reflect that failure in the status report:
#   - flip the row's Status from "Ready for handoff" to
#     "Handback transform failed" everywhere it is shown (Overview +
#     each language sheet), and
#   - record the mismatch detail in the Error Detail column of each
#     language sheet, widening that column so the message is readable.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$newStatus = "Handback transform failed"

$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

# Error Detail (column P) messages for the 5e890ed8 row, per language.
$wsZhCn.Range("P3").Value = "Handback file name: yyhf4t35.vza is different with handoff file name: 5e890ed8-fb68-466f-8fef-8f27989b91e8.0f2586832efd74b2fce3575b9e165489645d3de5.zh-cn."
$wsDeDe.Range("P3").Value = "Handback file name: yyhf4t35.vza is different with handoff file name: 5e890ed8-fb68-466f-8fef-8f27989b91e8.0f2586832efd74b2fce3575b9e165489645d3de5.de-de."

# Widen column P (Error Detail) on both language sheets to fit the new
# message text. The ColumnWidth COM property is offset from the raw
# OOXML column width (in "characters") by 5/6 of a character, so ask for
# 40 - 5/6 to land on an on-disk width of exactly 40.
$wsZhCn.Columns.Item(16).ColumnWidth = 39.166666666666664
$wsDeDe.Columns.Item(16).ColumnWidth = 39.166666666666664
